$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 73.57143000000001
$ws.Range("I4").Value = 73.57143000000001
$ws.Range("K4").Value = 73.57143000000001
$ws.Range("M4").Value = 40.42856999999999
$ws.Range("H62").Value = 3741.1667
$ws.Range("I62").Value = 1733
$ws.Range("K62").Value = 1733
$ws.Range("M62").Value = -1109
$ws.Range("H65").Value = 3741.1667
$ws.Range("I65").Value = 1733
$ws.Range("K65").Value = 8665
$ws.Range("M65").Value = -5545
$ws.Range("H92").Value = 1261.381
$ws.Range("I92").Value = 1249.0625
$ws.Range("K92").Value = 1249.0625
$ws.Range("M92").Value = -1.0625
$ws.Range("H137").Value = 11507.1
$ws.Range("I137").Value = 11320.667
$ws.Range("J137").Value = 12386
$ws.Range("K137").Value = 33962.001
$ws.Range("L137").Value = 37158
$ws.Range("M137").Value = -31412.001
$ws.Range("N137").Value = -42258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6993.8
$ws.Range("I2").Value = 5505
$ws.Range("J2").Value = 7986.3335
$ws.Range("K2").Value = 5505
$ws.Range("L2").Value = 7986.3335
$ws.Range("M2").Value = -5392
$ws.Range("N2").Value = -8212.333500000001
$ws.Range("H54").Value = 26247
$ws.Range("J54").Value = 26247
$ws.Range("L54").Value = 26247
$ws.Range("N54").Value = -27785
$ws.Range("H74").Value = 5219.212
$ws.Range("I74").Value = 5039.8066
$ws.Range("K74").Value = 5039.8066
$ws.Range("M74").Value = -4165.8066
$ws.Range("H77").Value = 5219.212
$ws.Range("I77").Value = 5039.8066
$ws.Range("K77").Value = 25199.033
$ws.Range("M77").Value = -20831.033
$ws.Range("H97").Value = 1625.1428
$ws.Range("I97").Value = 1509.8889
$ws.Range("K97").Value = 1509.8889
$ws.Range("M97").Value = -1013.8889
$ws.Range("H116").Value = 6993.8
$ws.Range("I116").Value = 5505
$ws.Range("J116").Value = 7986.3335
$ws.Range("K116").Value = 5505
$ws.Range("L116").Value = 7986.3335
$ws.Range("M116").Value = -3211
$ws.Range("N116").Value = -12574.3335
$ws.Range("H132").Value = 38552.332
$ws.Range("I132").Value = 2472.6
$ws.Range("K132").Value = 7417.799999999999
$ws.Range("M132").Value = -4887.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6993.8
$ws.Range("I3").Value = 5505
$ws.Range("J3").Value = 7986.3335
$ws.Range("K3").Value = 5505
$ws.Range("L3").Value = 7986.3335
$ws.Range("M3").Value = -5391
$ws.Range("N3").Value = -8214.333500000001
$ws.Range("H94").Value = 1653.2759
$ws.Range("I94").Value = 1247.4445
$ws.Range("K94").Value = 1247.4445
$ws.Range("M94").Value = -796.4445000000001
$ws.Range("H107").Value = 3581.1
$ws.Range("I107").Value = 3423.6
$ws.Range("J107").Value = 4053.6
$ws.Range("K107").Value = 3423.6
$ws.Range("L107").Value = 4053.6
$ws.Range("M107").Value = -1503.6
$ws.Range("N107").Value = -7893.6
$ws.Range("H134").Value = 1014
$ws.Range("I134").Value = 1014
$ws.Range("K134").Value = 3042
$ws.Range("M134").Value = -507

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4291.488
$ws.Range("I31").Value = 2343.75
$ws.Range("K31").Value = 2343.75
$ws.Range("M31").Value = -2048.75
$ws.Range("H34").Value = 4291.488
$ws.Range("I34").Value = 2343.75
$ws.Range("K34").Value = 2343.75
$ws.Range("M34").Value = -2141.75
$ws.Range("H58").Value = 3234.9333
$ws.Range("I58").Value = 3614.7144
$ws.Range("J58").Value = 2902.625
$ws.Range("K58").Value = 3614.7144
$ws.Range("L58").Value = 2902.625
$ws.Range("M58").Value = -3411.7144
$ws.Range("N58").Value = -3308.625
$ws.Range("H80").Value = 49977
$ws.Range("J80").Value = 49977
$ws.Range("L80").Value = 49977
$ws.Range("N80").Value = -52223
$ws.Range("H83").Value = 49977
$ws.Range("J83").Value = 49977
$ws.Range("L83").Value = 149931
$ws.Range("N83").Value = -161163
$ws.Range("H122").Value = 699
$ws.Range("I122").Value = 699
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2097
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 353
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3675.4211
$ws.Range("I132").Value = 3298.7
$ws.Range("K132").Value = 9896.099999999999
$ws.Range("M132").Value = -7366.099999999999
$ws.Range("H134").Value = 3429
$ws.Range("I134").Value = 3062.0833
$ws.Range("J134").Value = 4058
$ws.Range("K134").Value = 9186.249899999999
$ws.Range("L134").Value = 12174
$ws.Range("M134").Value = -6651.249899999999
$ws.Range("N134").Value = -17244
$ws.Range("H136").Value = 3234.9333
$ws.Range("I136").Value = 3614.7144
$ws.Range("J136").Value = 2902.625
$ws.Range("K136").Value = 10844.1432
$ws.Range("L136").Value = 8707.875
$ws.Range("M136").Value = -8294.143199999999
$ws.Range("N136").Value = -13807.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 6500
$ws.Range("I42").Value = 3000
$ws.Range("K42").Value = 9000
$ws.Range("M42").Value = -8466
$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2189
$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4944
$ws.Range("H107").Value = 780.36365
$ws.Range("J107").Value = 1041.4286
$ws.Range("L107").Value = 3124.2858
$ws.Range("N107").Value = -6964.2858
$ws.Range("H112").Value = 80673.336
$ws.Range("I112").Value = 70000
$ws.Range("J112").Value = 86010
$ws.Range("K112").Value = 210000
$ws.Range("L112").Value = 258030
$ws.Range("M112").Value = -208892
$ws.Range("N112").Value = -260246
$ws.Range("H114").Value = 15255.833
$ws.Range("I114").Value = 707
$ws.Range("J114").Value = 88000
$ws.Range("K114").Value = 2121
$ws.Range("L114").Value = 264000
$ws.Range("M114").Value = 1133
$ws.Range("N114").Value = -270508
$ws.Range("H134").Value = 3128.6
$ws.Range("I134").Value = 1571.9166
$ws.Range("J134").Value = 9355.333000000001
$ws.Range("K134").Value = 4715.7498
$ws.Range("L134").Value = 28065.999
$ws.Range("M134").Value = 354.2502000000004
$ws.Range("N134").Value = -38205.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3361.75
$ws.Range("I80").Value = 2999.6667
$ws.Range("J80").Value = 3579
$ws.Range("K80").Value = 2999.6667
$ws.Range("L80").Value = 3579
$ws.Range("M80").Value = -2001.6667
$ws.Range("N80").Value = -5575
$ws.Range("H83").Value = 3361.75
$ws.Range("I83").Value = 2999.6667
$ws.Range("J83").Value = 3579
$ws.Range("K83").Value = 14998.3335
$ws.Range("L83").Value = 17895
$ws.Range("M83").Value = -10006.3335
$ws.Range("N83").Value = -27879
$ws.Range("H97").Value = 1086.7826
$ws.Range("I97").Value = 713.17645
$ws.Range("K97").Value = 713.17645
$ws.Range("M97").Value = -217.17645
$ws.Range("H122").Value = 2359.8
$ws.Range("I122").Value = 1574.75
$ws.Range("K122").Value = 4724.25
$ws.Range("M122").Value = -2274.25
$ws.Range("H132").Value = 2834.9412
$ws.Range("I132").Value = 2614.2222
$ws.Range("J132").Value = 3083.25
$ws.Range("K132").Value = 7842.6666
$ws.Range("L132").Value = 9249.75
$ws.Range("M132").Value = -5312.6666
$ws.Range("N132").Value = -14309.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5244.353
$ws.Range("I7").Value = 4116.1
$ws.Range("K7").Value = 4116.1
$ws.Range("M7").Value = -4004.1
$ws.Range("H16").Value = 1423.1936
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 1564.8572
$ws.Range("I22").Value = 1176.75
$ws.Range("J22").Value = 2411.6365
$ws.Range("K22").Value = 1176.75
$ws.Range("L22").Value = 2411.6365
$ws.Range("M22").Value = -881.75
$ws.Range("N22").Value = -3001.6365
$ws.Range("H27").Value = 1564.8572
$ws.Range("I27").Value = 1176.75
$ws.Range("J27").Value = 2411.6365
$ws.Range("K27").Value = 1176.75
$ws.Range("L27").Value = 2411.6365
$ws.Range("M27").Value = -1069.75
$ws.Range("N27").Value = -2625.6365
$ws.Range("H46").Value = 1636.1364
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 4498.75
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 4498.75
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -4874.75
$ws.Range("H55").Value = 1104.65
$ws.Range("I55").Value = 697
$ws.Range("K55").Value = 697
$ws.Range("M55").Value = -524
$ws.Range("H126").Value = 5244.353
$ws.Range("I126").Value = 4116.1
$ws.Range("K126").Value = 12348.3
$ws.Range("M126").Value = -9878.300000000001
$ws.Range("H132").Value = 7592.375
$ws.Range("I132").Value = 1616.6923
$ws.Range("K132").Value = 4850.0769
$ws.Range("M132").Value = -2320.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1012.6667
$ws.Range("I132").Value = 864.5161000000001
$ws.Range("J132").Value = 1586.75
$ws.Range("K132").Value = 2593.5483
$ws.Range("L132").Value = 4760.25
$ws.Range("M132").Value = -63.54830000000038
$ws.Range("N132").Value = -9820.25
$ws.Range("H136").Value = 8306.32
$ws.Range("I136").Value = 7345.0312
$ws.Range("K136").Value = 22035.0936
$ws.Range("M136").Value = -19485.0936
